# Generate Report for Handback
# Updates the "Correspond Handoff Datetime" (col E) and "Correspond Handback
# DateTime" (col H) values for the first data row (row 2) on the per-locale
# report sheets, reflecting a fresh handback cycle.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E2").Value = "2016-03-20 04:45:23"
$zhcn.Range("H2").Value = "2016-03-20 04:46:29"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E2").Value = "2016-03-20 04:45:32"
$dede.Range("H2").Value = "2016-03-20 04:46:44"
